# Insert a new data row at sheet row 70 (pushing the existing rows 70-193
# down to 71-194, exactly as the target diff shows) and populate it with
# the new "Brócoli" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 70..193 down by one row.
$ws.Rows("70:70").Insert()

# Populate the newly inserted row 70 with its values.
$ws.Cells.Item(70, 1).Value  = 4
$ws.Cells.Item(70, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(70, 3).Value  = "Los Lagos"
$ws.Cells.Item(70, 4).Value  = 44469
$ws.Cells.Item(70, 5).Value  = 10
$ws.Cells.Item(70, 6).Value  = 100112023
$ws.Cells.Item(70, 7).Value  = "Brócoli"
$ws.Cells.Item(70, 8).Value  = "Sin especificar"
$ws.Cells.Item(70, 9).Value  = "Primera"
$ws.Cells.Item(70, 10).Value = 500
$ws.Cells.Item(70, 11).Value = 1300
$ws.Cells.Item(70, 12).Value = 1300
$ws.Cells.Item(70, 13).Value = 1300
$ws.Cells.Item(70, 14).Value = "$/unidad"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 1300
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(70, 18).Value = "Hortaliza"
